$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.735.80'
$ws.Range("E2").Value = '  -2.28%  '

$ws.Range("D3").Value = '3.207.11'
$ws.Range("E3").Value = '  -3.14%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.53%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.39'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.36%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = '3.207.84'
$ws.Range("E8").Value = '  -3.01%  '

$ws.Range("E9").Value = '  -3.03%  '

$ws.Range("E10").Value = '  -3.47%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.35'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.46%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.456'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.33%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000238'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.65%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.58'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.49%  '

$ws.Range("D15").Value = '3.736.20'
$ws.Range("E15").Value = '  -3.06%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.120'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.06%  '

$ws.Range("D17").Value = '3.204.49'
$ws.Range("E17").Value = '  -3.18%  '

$ws.Range("D18").Value = '62.811.49'
$ws.Range("E18").Value = '  -2.28%  '

$ws.Range("E19").Value = '  -3.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '465.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.48%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.90'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.45%  '

$ws.Range("E22").Value = '  -4.35%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.69'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.96%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.35'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.95%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.41'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.73%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.01%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.68'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.80%  '

$ws.Range("E28").Value = '  -0.18%  '

$ws.Range("E29").Value = '  -5.32%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.93'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.61%  '

$ws.Range("E31").Value = '  -4.23%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.65'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.89%  '

$ws.Range("E33").Value = '  -5.64%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.43'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.17%  '

$ws.Range("E35").Value = '  -5.13%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.86'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.99%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.80'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.81%  '

$ws.Range("D38").Value = '0.0₃0700'
$ws.Range("E38").Value = '  -5.78%  '

$ws.Range("E39").Value = '  -2.28%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '420.19'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.51%  '

$ws.Range("D41").Value = '3.014.62'
$ws.Range("E41").Value = '  -0.25%  '

$ws.Range("E42").Value = '  +4.47%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.11'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.63%  '

$ws.Range("E44").Value = '  -6.43%  '

$ws.Range("E45").Value = '  -6.44%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.14'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.16%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '35.74'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.10%  '

$ws.Range("E49").Value = '  -2.98%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '125.23'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.20%  '

$ws.Range("E51").Value = '  -3.04%  '
